$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 389
$ws.Range("D389").Value = 44641
$ws.Range("H389").Value = "Cuatro cascos rojo"
$ws.Range("I389").Value = "Primera"
$ws.Range("J389").Value = 300
$ws.Range("K389").Value = 13000
$ws.Range("L389").Value = 13000
$ws.Range("M389").Value = 13000
$ws.Range("N389").Value = "`$/caja 15 kilos"
$ws.Range("O389").Value = "Región del Maule"
$ws.Range("P389").Value = 867
$ws.Range("Q389").Value = 15

# Row 390
$ws.Range("D390").Value = 44641
$ws.Range("H390").Value = "Cuatro cascos verde"
$ws.Range("I390").Value = "Primera"
$ws.Range("J390").Value = 300
$ws.Range("K390").Value = 7000
$ws.Range("L390").Value = 7000
$ws.Range("M390").Value = 7000
$ws.Range("N390").Value = "`$/caja 15 kilos"
$ws.Range("O390").Value = "Región del Maule"
$ws.Range("P390").Value = 467
$ws.Range("Q390").Value = 15

# Row 391
$ws.Range("D391").Value = 44421
$ws.Range("H391").Value = "Zafiro rojo"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 300
$ws.Range("K391").Value = 40000
$ws.Range("L391").Value = 40000
$ws.Range("M391").Value = 40000
$ws.Range("N391").Value = "`$/caja 15 kilos"
$ws.Range("O391").Value = "Región de Arica y Parinacota"
$ws.Range("P391").Value = 2667
$ws.Range("Q391").Value = 15

# Row 392
$ws.Range("D392").Value = 44421
$ws.Range("H392").Value = "Zafiro verde"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 300
$ws.Range("K392").Value = 38000
$ws.Range("L392").Value = 38000
$ws.Range("M392").Value = 38000
$ws.Range("N392").Value = "`$/caja 15 kilos"
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 2533
$ws.Range("Q392").Value = 15

# Row 393
$ws.Range("D393").Value = 44329
$ws.Range("H393").Value = "Zafiro rojo"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 300
$ws.Range("K393").Value = 20000
$ws.Range("L393").Value = 20000
$ws.Range("M393").Value = 20000
$ws.Range("N393").Value = "`$/caja 15 kilos"
$ws.Range("O393").Value = "Región de Arica y Parinacota"
$ws.Range("P393").Value = 1333
$ws.Range("Q393").Value = 15

# Row 394
$ws.Range("D394").Value = 44329
$ws.Range("H394").Value = "Zafiro verde"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 300
$ws.Range("K394").Value = 15000
$ws.Range("L394").Value = 15000
$ws.Range("M394").Value = 15000
$ws.Range("N394").Value = "`$/caja 15 kilos"
$ws.Range("O394").Value = "Región de Arica y Parinacota"
$ws.Range("P394").Value = 1000
$ws.Range("Q394").Value = 15

# Row 395
$ws.Range("D395").Value = 44637
$ws.Range("H395").Value = "Cuatro cascos rojo"
$ws.Range("I395").Value = "Primera"
$ws.Range("J395").Value = 200
$ws.Range("K395").Value = 15000
$ws.Range("L395").Value = 15000
$ws.Range("M395").Value = 15000
$ws.Range("N395").Value = "`$/caja 15 kilos"
$ws.Range("O395").Value = "Región del Maule"
$ws.Range("P395").Value = 1000
$ws.Range("Q395").Value = 15

# Row 396
$ws.Range("D396").Value = 44637
$ws.Range("H396").Value = "Cuatro cascos verde"
$ws.Range("I396").Value = "Primera"
$ws.Range("J396").Value = 300
$ws.Range("K396").Value = 8000
$ws.Range("L396").Value = 8000
$ws.Range("M396").Value = 8000
$ws.Range("N396").Value = "`$/caja 15 kilos"
$ws.Range("O396").Value = "Región del Maule"
$ws.Range("P396").Value = 533
$ws.Range("Q396").Value = 15

# Row 397
$ws.Range("D397").Value = 44637
$ws.Range("H397").Value = "Zafiro rojo"
$ws.Range("I397").Value = "Primera"
$ws.Range("J397").Value = 200
$ws.Range("K397").Value = 18000
$ws.Range("L397").Value = 18000
$ws.Range("M397").Value = 18000
$ws.Range("N397").Value = "`$/caja 15 kilos"
$ws.Range("O397").Value = "Región de Arica y Parinacota"
$ws.Range("P397").Value = 1200
$ws.Range("Q397").Value = 15

# Row 398
$ws.Range("D398").Value = 44208
$ws.Range("H398").Value = "Cuatro cascos verde"
$ws.Range("I398").Value = "Primera"
$ws.Range("J398").Value = 200
$ws.Range("K398").Value = 10000
$ws.Range("L398").Value = 10000
$ws.Range("M398").Value = 10000
$ws.Range("N398").Value = "`$/caja 15 kilos"
$ws.Range("O398").Value = "Región del Maule"
$ws.Range("P398").Value = 667
$ws.Range("Q398").Value = 15

# Row 399
$ws.Range("D399").Value = 44445
$ws.Range("H399").Value = "Zafiro rojo"
$ws.Range("I399").Value = "Primera"
$ws.Range("J399").Value = 200
$ws.Range("K399").Value = 42000
$ws.Range("L399").Value = 42000
$ws.Range("M399").Value = 42000
$ws.Range("N399").Value = "`$/caja 15 kilos"
$ws.Range("O399").Value = "Región de Arica y Parinacota"
$ws.Range("P399").Value = 2800
$ws.Range("Q399").Value = 15

# Row 400
$ws.Range("D400").Value = 44445
$ws.Range("H400").Value = "Zafiro verde"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 200
$ws.Range("K400").Value = 38000
$ws.Range("L400").Value = 38000
$ws.Range("M400").Value = 38000
$ws.Range("N400").Value = "`$/caja 15 kilos"
$ws.Range("O400").Value = "Región de Arica y Parinacota"
$ws.Range("P400").Value = 2533
$ws.Range("Q400").Value = 15

# Row 401
$ws.Range("D401").Value = 44524
$ws.Range("H401").Value = "Cuatro cascos verde"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 200
$ws.Range("K401").Value = 15000
$ws.Range("L401").Value = 15000
$ws.Range("M401").Value = 15000
$ws.Range("N401").Value = "`$/caja 15 kilos"
$ws.Range("O401").Value = "Región del Maule"
$ws.Range("P401").Value = 1000
$ws.Range("Q401").Value = 15

# Row 402
$ws.Range("D402").Value = 44355
$ws.Range("H402").Value = "Morrón rojo"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 300
$ws.Range("K402").Value = 15000
$ws.Range("L402").Value = 15000
$ws.Range("M402").Value = 15000
$ws.Range("N402").Value = "`$/caja 18 kilos"
$ws.Range("O402").Value = "Provincia del Elquí"
$ws.Range("P402").Value = 833
$ws.Range("Q402").Value = 18

# Row 403
$ws.Range("D403").Value = 44355
$ws.Range("H403").Value = "Zafiro rojo"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 300
$ws.Range("K403").Value = 15000
$ws.Range("L403").Value = 15000
$ws.Range("M403").Value = 15000
$ws.Range("N403").Value = "`$/caja 15 kilos"
$ws.Range("O403").Value = "Región de Arica y Parinacota"
$ws.Range("P403").Value = 1000
$ws.Range("Q403").Value = 15

# Row 404
$ws.Range("D404").Value = 44355
$ws.Range("H404").Value = "Zafiro verde"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 300
$ws.Range("K404").Value = 11000
$ws.Range("L404").Value = 11000
$ws.Range("M404").Value = 11000
$ws.Range("N404").Value = "`$/caja 15 kilos"
$ws.Range("O404").Value = "Región de Arica y Parinacota"
$ws.Range("P404").Value = 733
$ws.Range("Q404").Value = 15

# Row 405
$ws.Range("D405").Value = 44530
$ws.Range("H405").Value = "Cuatro cascos verde"
$ws.Range("I405").Value = "Primera"
$ws.Range("J405").Value = 200
$ws.Range("K405").Value = 15000
$ws.Range("L405").Value = 15000
$ws.Range("M405").Value = 15000
$ws.Range("N405").Value = "`$/caja 15 kilos"
$ws.Range("O405").Value = "Región del Maule"
$ws.Range("P405").Value = 1000
$ws.Range("Q405").Value = 15

# Row 406
$ws.Range("D406").Value = 44483
$ws.Range("H406").Value = "Zafiro rojo"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 200
$ws.Range("K406").Value = 42000
$ws.Range("L406").Value = 42000
$ws.Range("M406").Value = 42000
$ws.Range("N406").Value = "`$/caja 15 kilos"
$ws.Range("O406").Value = "Región de Arica y Parinacota"
$ws.Range("P406").Value = 2800
$ws.Range("Q406").Value = 15

# Row 407
$ws.Range("D407").Value = 44483
$ws.Range("H407").Value = "Zafiro verde"
$ws.Range("I407").Value = "Primera"
$ws.Range("J407").Value = 200
$ws.Range("K407").Value = 40000
$ws.Range("L407").Value = 40000
$ws.Range("M407").Value = 40000
$ws.Range("N407").Value = "`$/caja 15 kilos"
$ws.Range("O407").Value = "Región de Arica y Parinacota"
$ws.Range("P407").Value = 2667
$ws.Range("Q407").Value = 15

# Row 408
$ws.Range("D408").Value = 44617
$ws.Range("H408").Value = "Cuatro cascos verde"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 300
$ws.Range("K408").Value = 7000
$ws.Range("L408").Value = 7000
$ws.Range("M408").Value = 7000
$ws.Range("N408").Value = "`$/caja 15 kilos"
$ws.Range("O408").Value = "Región del Maule"
$ws.Range("P408").Value = 467
$ws.Range("Q408").Value = 15

# Row 409
$ws.Range("D409").Value = 44557
$ws.Range("H409").Value = "Cuatro cascos verde"
$ws.Range("I409").Value = "Primera"
$ws.Range("J409").Value = 300
$ws.Range("K409").Value = 10000
$ws.Range("L409").Value = 10000
$ws.Range("M409").Value = 10000
$ws.Range("N409").Value = "`$/caja 15 kilos"
$ws.Range("O409").Value = "Región del Maule"
$ws.Range("P409").Value = 667
$ws.Range("Q409").Value = 15

# Row 410
$ws.Range("D410").Value = 44489
$ws.Range("H410").Value = "Zafiro rojo"
$ws.Range("I410").Value = "Primera"
$ws.Range("J410").Value = 150
$ws.Range("K410").Value = 43000
$ws.Range("L410").Value = 43000
$ws.Range("M410").Value = 43000
$ws.Range("N410").Value = "`$/caja 15 kilos"
$ws.Range("O410").Value = "Región de Arica y Parinacota"
$ws.Range("P410").Value = 2867
$ws.Range("Q410").Value = 15

# Row 411
$ws.Range("D411").Value = 44264
$ws.Range("H411").Value = "Cuatro cascos rojo"
$ws.Range("I411").Value = "Primera"
$ws.Range("J411").Value = 200
$ws.Range("K411").Value = 12000
$ws.Range("L411").Value = 12000
$ws.Range("M411").Value = 12000
$ws.Range("N411").Value = "`$/caja 15 kilos"
$ws.Range("O411").Value = "Región del Maule"
$ws.Range("P411").Value = 800
$ws.Range("Q411").Value = 15

# Row 412
$ws.Range("D412").Value = 44264
$ws.Range("H412").Value = "Cuatro cascos verde"
$ws.Range("I412").Value = "Primera"
$ws.Range("J412").Value = 200
$ws.Range("K412").Value = 6000
$ws.Range("L412").Value = 6000
$ws.Range("M412").Value = 6000
$ws.Range("N412").Value = "`$/caja 15 kilos"
$ws.Range("O412").Value = "Región del Maule"
$ws.Range("P412").Value = 400
$ws.Range("Q412").Value = 15

# Row 413
$ws.Range("D413").Value = 44396
$ws.Range("H413").Value = "Cuatro cascos rojo"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 300
$ws.Range("K413").Value = 17000
$ws.Range("L413").Value = 17000
$ws.Range("M413").Value = 17000
$ws.Range("N413").Value = "`$/caja 15 kilos"
$ws.Range("O413").Value = "Provincia del Elquí"
$ws.Range("P413").Value = 1133
$ws.Range("Q413").Value = 15

# Row 414
$ws.Range("D414").Value = 44396
$ws.Range("H414").Value = "Cuatro cascos rojo"
$ws.Range("I414").Value = "Segunda"
$ws.Range("J414").Value = 200
$ws.Range("K414").Value = 14000
$ws.Range("L414").Value = 14000
$ws.Range("M414").Value = 14000
$ws.Range("N414").Value = "`$/caja 15 kilos"
$ws.Range("O414").Value = "Provincia del Elquí"
$ws.Range("P414").Value = 933
$ws.Range("Q414").Value = 15

# Row 415
$ws.Range("D415").Value = 44396
$ws.Range("H415").Value = "Cuatro cascos verde"
$ws.Range("I415").Value = "Primera"
$ws.Range("J415").Value = 300
$ws.Range("K415").Value = 20000
$ws.Range("L415").Value = 20000
$ws.Range("M415").Value = 20000
$ws.Range("N415").Value = "`$/caja 15 kilos"
$ws.Range("O415").Value = "Provincia del Elquí"
$ws.Range("P415").Value = 1333
$ws.Range("Q415").Value = 15

# Row 416
$ws.Range("D416").Value = 44232
$ws.Range("H416").Value = "Cuatro cascos verde"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 300
$ws.Range("K416").Value = 7000
$ws.Range("L416").Value = 7000
$ws.Range("M416").Value = 7000
$ws.Range("N416").Value = "`$/caja 15 kilos"
$ws.Range("O416").Value = "Región del Maule"
$ws.Range("P416").Value = 467
$ws.Range("Q416").Value = 15

# Row 417
$ws.Range("D417").Value = 44279
$ws.Range("H417").Value = "Cuatro cascos verde"
$ws.Range("I417").Value = "Primera"
$ws.Range("J417").Value = 200
$ws.Range("K417").Value = 7000
$ws.Range("L417").Value = 7000
$ws.Range("M417").Value = 7000
$ws.Range("N417").Value = "`$/caja 15 kilos"
$ws.Range("O417").Value = "Región del Maule"
$ws.Range("P417").Value = 467
$ws.Range("Q417").Value = 15

# Row 418
$ws.Range("D418").Value = 44330
$ws.Range("H418").Value = "Zafiro rojo"
$ws.Range("I418").Value = "Primera"
$ws.Range("J418").Value = 200
$ws.Range("K418").Value = 23000
$ws.Range("L418").Value = 23000
$ws.Range("M418").Value = 23000
$ws.Range("N418").Value = "`$/caja 15 kilos"
$ws.Range("O418").Value = "Región de Arica y Parinacota"
$ws.Range("P418").Value = 1533
$ws.Range("Q418").Value = 15

# Row 419
$ws.Range("D419").Value = 44330
$ws.Range("H419").Value = "Zafiro verde"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 300
$ws.Range("K419").Value = 15000
$ws.Range("L419").Value = 15000
$ws.Range("M419").Value = 15000
$ws.Range("N419").Value = "`$/caja 15 kilos"
$ws.Range("O419").Value = "Región de Arica y Parinacota"
$ws.Range("P419").Value = 1000
$ws.Range("Q419").Value = 15

# Row 420
$ws.Range("D420").Value = 44504
$ws.Range("H420").Value = "Zafiro rojo"
$ws.Range("I420").Value = "Primera"
$ws.Range("J420").Value = 200
$ws.Range("K420").Value = 43000
$ws.Range("L420").Value = 43000
$ws.Range("M420").Value = 43000
$ws.Range("N420").Value = "`$/caja 15 kilos"
$ws.Range("O420").Value = "Región de Arica y Parinacota"
$ws.Range("P420").Value = 2867
$ws.Range("Q420").Value = 15

# Row 421
$ws.Range("D421").Value = 44504
$ws.Range("H421").Value = "Zafiro verde"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 300
$ws.Range("K421").Value = 35000
$ws.Range("L421").Value = 35000
$ws.Range("M421").Value = 35000
$ws.Range("N421").Value = "`$/caja 15 kilos"
$ws.Range("O421").Value = "Región de Arica y Parinacota"
$ws.Range("P421").Value = 2333
$ws.Range("Q421").Value = 15

# Row 422
$ws.Range("D422").Value = 44572
$ws.Range("H422").Value = "Cuatro cascos verde"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 200
$ws.Range("K422").Value = 9000
$ws.Range("L422").Value = 9000
$ws.Range("M422").Value = 9000
$ws.Range("N422").Value = "`$/caja 15 kilos"
$ws.Range("O422").Value = "Región del Maule"
$ws.Range("P422").Value = 600
$ws.Range("Q422").Value = 15

# Row 423
$ws.Range("D423").Value = 44257
$ws.Range("H423").Value = "Cuatro cascos rojo"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 150
$ws.Range("K423").Value = 9000
$ws.Range("L423").Value = 9000
$ws.Range("M423").Value = 9000
$ws.Range("N423").Value = "`$/caja 15 kilos"
$ws.Range("O423").Value = "Región del Maule"
$ws.Range("P423").Value = 600
$ws.Range("Q423").Value = 15

# Row 424
$ws.Range("D424").Value = 44257
$ws.Range("H424").Value = "Cuatro cascos verde"
$ws.Range("I424").Value = "Primera"
$ws.Range("J424").Value = 300
$ws.Range("K424").Value = 5000
$ws.Range("L424").Value = 5000
$ws.Range("M424").Value = 5000
$ws.Range("N424").Value = "`$/caja 15 kilos"
$ws.Range("O424").Value = "Región del Maule"
$ws.Range("P424").Value = 333
$ws.Range("Q424").Value = 15

# Row 425
$ws.Range("D425").Value = 44301
$ws.Range("H425").Value = "Zafiro rojo"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 200
$ws.Range("K425").Value = 16000
$ws.Range("L425").Value = 16000
$ws.Range("M425").Value = 16000
$ws.Range("N425").Value = "`$/caja 15 kilos"
$ws.Range("O425").Value = "Región de Arica y Parinacota"
$ws.Range("P425").Value = 1067
$ws.Range("Q425").Value = 15

# Row 426
$ws.Range("D426").Value = 44301
$ws.Range("H426").Value = "Zafiro verde"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 300
$ws.Range("K426").Value = 10000
$ws.Range("L426").Value = 10000
$ws.Range("M426").Value = 10000
$ws.Range("N426").Value = "`$/caja 15 kilos"
$ws.Range("O426").Value = "Región de Arica y Parinacota"
$ws.Range("P426").Value = 667
$ws.Range("Q426").Value = 15

# Row 427
$ws.Range("D427").Value = 44370
$ws.Range("H427").Value = "Zafiro rojo"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 400
$ws.Range("K427").Value = 15000
$ws.Range("L427").Value = 15000
$ws.Range("M427").Value = 15000
$ws.Range("N427").Value = "`$/caja 15 kilos"
$ws.Range("O427").Value = "Región de Arica y Parinacota"
$ws.Range("P427").Value = 1000
$ws.Range("Q427").Value = 15

# Row 428
$ws.Range("D428").Value = 44370
$ws.Range("H428").Value = "Zafiro verde"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 400
$ws.Range("K428").Value = 11000
$ws.Range("L428").Value = 11000
$ws.Range("M428").Value = 11000
$ws.Range("N428").Value = "`$/caja 15 kilos"
$ws.Range("O428").Value = "Región de Arica y Parinacota"
$ws.Range("P428").Value = 733
$ws.Range("Q428").Value = 15

# Row 429
$ws.Range("D429").Value = 44487
$ws.Range("H429").Value = "Zafiro rojo"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 200
$ws.Range("K429").Value = 43000
$ws.Range("L429").Value = 43000
$ws.Range("M429").Value = 43000
$ws.Range("N429").Value = "`$/caja 15 kilos"
$ws.Range("O429").Value = "Región de Arica y Parinacota"
$ws.Range("P429").Value = 2867
$ws.Range("Q429").Value = 15

# Row 430
$ws.Range("D430").Value = 44487
$ws.Range("H430").Value = "Zafiro verde"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 300
$ws.Range("K430").Value = 40000
$ws.Range("L430").Value = 40000
$ws.Range("M430").Value = 40000
$ws.Range("N430").Value = "`$/caja 15 kilos"
$ws.Range("O430").Value = "Región de Arica y Parinacota"
$ws.Range("P430").Value = 2667
$ws.Range("Q430").Value = 15

# Row 431
$ws.Range("D431").Value = 44174
$ws.Range("H431").Value = "Cuatro cascos verde"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 250
$ws.Range("K431").Value = 15000
$ws.Range("L431").Value = 15000
$ws.Range("M431").Value = 15000
$ws.Range("N431").Value = "`$/caja 15 kilos"
$ws.Range("O431").Value = "Región del Maule"
$ws.Range("P431").Value = 1000
$ws.Range("Q431").Value = 15

# Row 432
$ws.Range("D432").Value = 44200
$ws.Range("H432").Value = "Cuatro cascos verde"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 200
$ws.Range("K432").Value = 10000
$ws.Range("L432").Value = 12000
$ws.Range("M432").Value = 11000
$ws.Range("N432").Value = "`$/caja 15 kilos"
$ws.Range("O432").Value = "Región del Maule"
$ws.Range("P432").Value = 733
$ws.Range("Q432").Value = 15

# Row 433
$ws.Range("D433").Value = 44236
$ws.Range("H433").Value = "Cuatro cascos rojo"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 200
$ws.Range("K433").Value = 13000
$ws.Range("L433").Value = 13000
$ws.Range("M433").Value = 13000
$ws.Range("N433").Value = "`$/caja 15 kilos"
$ws.Range("O433").Value = "Región del Maule"
$ws.Range("P433").Value = 867
$ws.Range("Q433").Value = 15

# Row 434
$ws.Range("D434").Value = 44236
$ws.Range("H434").Value = "Cuatro cascos verde"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 400
$ws.Range("K434").Value = 6000
$ws.Range("L434").Value = 6000
$ws.Range("M434").Value = 6000
$ws.Range("N434").Value = "`$/caja 15 kilos"
$ws.Range("O434").Value = "Región del Maule"
$ws.Range("P434").Value = 400
$ws.Range("Q434").Value = 15

# Row 435
$ws.Range("D435").Value = 44221
$ws.Range("H435").Value = "Cuatro cascos rojo"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 200
$ws.Range("K435").Value = 13000
$ws.Range("L435").Value = 13000
$ws.Range("M435").Value = 13000
$ws.Range("N435").Value = "`$/caja 15 kilos"
$ws.Range("O435").Value = "Región del Maule"
$ws.Range("P435").Value = 867
$ws.Range("Q435").Value = 15

# Row 436
$ws.Range("D436").Value = 44221
$ws.Range("H436").Value = "Cuatro cascos verde"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 300
$ws.Range("K436").Value = 9000
$ws.Range("L436").Value = 9000
$ws.Range("M436").Value = 9000
$ws.Range("N436").Value = "`$/caja 15 kilos"
$ws.Range("O436").Value = "Región del Maule"
$ws.Range("P436").Value = 600
$ws.Range("Q436").Value = 15

# Row 437
$ws.Range("D437").Value = 44413
$ws.Range("H437").Value = "Zafiro rojo"
$ws.Range("I437").Value = "Primera"
$ws.Range("J437").Value = 300
$ws.Range("K437").Value = 25000
$ws.Range("L437").Value = 25000
$ws.Range("M437").Value = 25000
$ws.Range("N437").Value = "`$/caja 15 kilos"
$ws.Range("O437").Value = "Región de Arica y Parinacota"
$ws.Range("P437").Value = 1667
$ws.Range("Q437").Value = 15

# Row 438
$ws.Range("D438").Value = 44413
$ws.Range("H438").Value = "Zafiro verde"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 300
$ws.Range("K438").Value = 20000
$ws.Range("L438").Value = 20000
$ws.Range("M438").Value = 20000
$ws.Range("N438").Value = "`$/caja 15 kilos"
$ws.Range("O438").Value = "Región de Arica y Parinacota"
$ws.Range("P438").Value = 1333
$ws.Range("Q438").Value = 15

# Row 439
$ws.Range("D439").Value = 44272
$ws.Range("H439").Value = "Cuatro cascos rojo"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 200
$ws.Range("K439").Value = 10000
$ws.Range("L439").Value = 10000
$ws.Range("M439").Value = 10000
$ws.Range("N439").Value = "`$/caja 15 kilos"
$ws.Range("O439").Value = "Región del Maule"
$ws.Range("P439").Value = 667
$ws.Range("Q439").Value = 15

# Row 440
$ws.Range("D440").Value = 44272
$ws.Range("H440").Value = "Cuatro cascos verde"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 200
$ws.Range("K440").Value = 6000
$ws.Range("L440").Value = 6000
$ws.Range("M440").Value = 6000
$ws.Range("N440").Value = "`$/caja 15 kilos"
$ws.Range("O440").Value = "Región del Maule"
$ws.Range("P440").Value = 400
$ws.Range("Q440").Value = 15

# Row 441
$ws.Range("D441").Value = 44229
$ws.Range("H441").Value = "Cuatro cascos verde"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 400
$ws.Range("K441").Value = 7000
$ws.Range("L441").Value = 7000
$ws.Range("M441").Value = 7000
$ws.Range("N441").Value = "`$/caja 15 kilos"
$ws.Range("O441").Value = "Región del Maule"
$ws.Range("P441").Value = 467
$ws.Range("Q441").Value = 15

# Row 442
$ws.Range("D442").Value = 44214
$ws.Range("H442").Value = "Cuatro cascos verde"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 300
$ws.Range("K442").Value = 8000
$ws.Range("L442").Value = 9000
$ws.Range("M442").Value = 8500
$ws.Range("N442").Value = "`$/caja 15 kilos"
$ws.Range("O442").Value = "Región del Maule"
$ws.Range("P442").Value = 567
$ws.Range("Q442").Value = 15

# Row 443
$ws.Range("D443").Value = 44299
$ws.Range("H443").Value = "Cuatro cascos rojo"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 200
$ws.Range("K443").Value = 12000
$ws.Range("L443").Value = 12000
$ws.Range("M443").Value = 12000
$ws.Range("N443").Value = "`$/caja 15 kilos"
$ws.Range("O443").Value = "Región del Maule"
$ws.Range("P443").Value = 800
$ws.Range("Q443").Value = 15

# Row 444
$ws.Range("D444").Value = 44299
$ws.Range("H444").Value = "Cuatro cascos verde"
$ws.Range("I444").Value = "Primera"
$ws.Range("J444").Value = 300
$ws.Range("K444").Value = 8000
$ws.Range("L444").Value = 8000
$ws.Range("M444").Value = 8000
$ws.Range("N444").Value = "`$/caja 15 kilos"
$ws.Range("O444").Value = "Región del Maule"
$ws.Range("P444").Value = 533
$ws.Range("Q444").Value = 15

# Row 445
$ws.Range("D445").Value = 44610
$ws.Range("H445").Value = "Cuatro cascos rojo"
$ws.Range("I445").Value = "Primera"
$ws.Range("J445").Value = 200
$ws.Range("K445").Value = 10000
$ws.Range("L445").Value = 10000
$ws.Range("M445").Value = 10000
$ws.Range("N445").Value = "`$/caja 15 kilos"
$ws.Range("O445").Value = "Región del Maule"
$ws.Range("P445").Value = 667
$ws.Range("Q445").Value = 15

# Row 446
$ws.Range("D446").Value = 44610
$ws.Range("H446").Value = "Cuatro cascos verde"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 300
$ws.Range("K446").Value = 6000
$ws.Range("L446").Value = 6000
$ws.Range("M446").Value = 6000
$ws.Range("N446").Value = "`$/caja 15 kilos"
$ws.Range("O446").Value = "Región del Maule"
$ws.Range("P446").Value = 400
$ws.Range("Q446").Value = 15

# Row 447
$ws.Range("D447").Value = 44312
$ws.Range("H447").Value = "Cuatro cascos rojo"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 300
$ws.Range("K447").Value = 9000
$ws.Range("L447").Value = 9000
$ws.Range("M447").Value = 9000
$ws.Range("N447").Value = "`$/caja 15 kilos"
$ws.Range("O447").Value = "Región del Maule"
$ws.Range("P447").Value = 600
$ws.Range("Q447").Value = 15

# Row 448
$ws.Range("D448").Value = 44312
$ws.Range("H448").Value = "Cuatro cascos verde"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 400
$ws.Range("K448").Value = 6000
$ws.Range("L448").Value = 6000
$ws.Range("M448").Value = 6000
$ws.Range("N448").Value = "`$/caja 15 kilos"
$ws.Range("O448").Value = "Región del Maule"
$ws.Range("P448").Value = 400
$ws.Range("Q448").Value = 15

# Row 449
$ws.Range("D449").Value = 44399
$ws.Range("H449").Value = "Zafiro rojo"
$ws.Range("I449").Value = "Primera"
$ws.Range("J449").Value = 300
$ws.Range("K449").Value = 17000
$ws.Range("L449").Value = 17000
$ws.Range("M449").Value = 17000
$ws.Range("N449").Value = "`$/caja 15 kilos"
$ws.Range("O449").Value = "Región de Arica y Parinacota"
$ws.Range("P449").Value = 1133
$ws.Range("Q449").Value = 15

# Row 450
$ws.Range("D450").Value = 44399
$ws.Range("H450").Value = "Zafiro verde"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 300
$ws.Range("K450").Value = 15000
$ws.Range("L450").Value = 15000
$ws.Range("M450").Value = 15000
$ws.Range("N450").Value = "`$/caja 15 kilos"
$ws.Range("O450").Value = "Región de Arica y Parinacota"
$ws.Range("P450").Value = 1000
$ws.Range("Q450").Value = 15

# Row 451
$ws.Range("D451").Value = 44615
$ws.Range("H451").Value = "Cuatro cascos verde"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 300
$ws.Range("K451").Value = 6000
$ws.Range("L451").Value = 6000
$ws.Range("M451").Value = 6000
$ws.Range("N451").Value = "`$/caja 15 kilos"
$ws.Range("O451").Value = "Región del Maule"
$ws.Range("P451").Value = 400
$ws.Range("Q451").Value = 15

# Row 452
$ws.Range("D452").Value = 44522
$ws.Range("H452").Value = "Cuatro cascos verde"
$ws.Range("I452").Value = "Primera"
$ws.Range("J452").Value = 150
$ws.Range("K452").Value = 20000
$ws.Range("L452").Value = 20000
$ws.Range("M452").Value = 20000
$ws.Range("N452").Value = "`$/caja 15 kilos"
$ws.Range("O452").Value = "Región del Maule"
$ws.Range("P452").Value = 1333
$ws.Range("Q452").Value = 15

# Row 453
$ws.Range("D453").Value = 44543
$ws.Range("H453").Value = "Cuatro cascos verde"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 300
$ws.Range("K453").Value = 13000
$ws.Range("L453").Value = 13000
$ws.Range("M453").Value = 13000
$ws.Range("N453").Value = "`$/caja 15 kilos"
$ws.Range("O453").Value = "Región del Maule"
$ws.Range("P453").Value = 867
$ws.Range("Q453").Value = 15

# Row 454
$ws.Range("D454").Value = 44167
$ws.Range("H454").Value = "Cuatro cascos verde"
$ws.Range("I454").Value = "Primera"
$ws.Range("J454").Value = 200
$ws.Range("K454").Value = 15000
$ws.Range("L454").Value = 15000
$ws.Range("M454").Value = 15000
$ws.Range("N454").Value = "`$/caja 15 kilos"
$ws.Range("O454").Value = "Región del Maule"
$ws.Range("P454").Value = 1000
$ws.Range("Q454").Value = 15

# Row 455
$ws.Range("D455").Value = 44277
$ws.Range("H455").Value = "Cuatro cascos rojo"
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 200
$ws.Range("K455").Value = 10000
$ws.Range("L455").Value = 10000
$ws.Range("M455").Value = 10000
$ws.Range("N455").Value = "`$/caja 15 kilos"
$ws.Range("O455").Value = "Región del Maule"
$ws.Range("P455").Value = 667
$ws.Range("Q455").Value = 15

# Row 456
$ws.Range("D456").Value = 44277
$ws.Range("H456").Value = "Cuatro cascos verde"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 300
$ws.Range("K456").Value = 6500
$ws.Range("L456").Value = 6500
$ws.Range("M456").Value = 6500
$ws.Range("N456").Value = "`$/caja 15 kilos"
$ws.Range("O456").Value = "Región del Maule"
$ws.Range("P456").Value = 433
$ws.Range("Q456").Value = 15

# Row 457
$ws.Range("D457").Value = 44258
$ws.Range("H457").Value = "Cuatro cascos rojo"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 300
$ws.Range("K457").Value = 9000
$ws.Range("L457").Value = 9000
$ws.Range("M457").Value = 9000
$ws.Range("N457").Value = "`$/caja 15 kilos"
$ws.Range("O457").Value = "Región del Maule"
$ws.Range("P457").Value = 600
$ws.Range("Q457").Value = 15

# Row 458
$ws.Range("D458").Value = 44258
$ws.Range("H458").Value = "Cuatro cascos verde"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 300
$ws.Range("K458").Value = 6000
$ws.Range("L458").Value = 6000
$ws.Range("M458").Value = 6000
$ws.Range("N458").Value = "`$/caja 15 kilos"
$ws.Range("O458").Value = "Región del Maule"
$ws.Range("P458").Value = 400
$ws.Range("Q458").Value = 15

# Row 459
$ws.Range("D459").Value = 44390
$ws.Range("H459").Value = "Zafiro verde"
$ws.Range("I459").Value = "Primera"
$ws.Range("J459").Value = 300
$ws.Range("K459").Value = 14000
$ws.Range("L459").Value = 14000
$ws.Range("M459").Value = 14000
$ws.Range("N459").Value = "`$/caja 15 kilos"
$ws.Range("O459").Value = "Región de Arica y Parinacota"
$ws.Range("P459").Value = 933
$ws.Range("Q459").Value = 15

# Row 460
$ws.Range("D460").Value = 44349
$ws.Range("H460").Value = "Zafiro rojo"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 300
$ws.Range("K460").Value = 20000
$ws.Range("L460").Value = 20000
$ws.Range("M460").Value = 20000
$ws.Range("N460").Value = "`$/caja 15 kilos"
$ws.Range("O460").Value = "Región de Arica y Parinacota"
$ws.Range("P460").Value = 1333
$ws.Range("Q460").Value = 15

# Row 461
$ws.Range("D461").Value = 44349
$ws.Range("H461").Value = "Zafiro verde"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 400
$ws.Range("K461").Value = 11000
$ws.Range("L461").Value = 11000
$ws.Range("M461").Value = 11000
$ws.Range("N461").Value = "`$/caja 15 kilos"
$ws.Range("O461").Value = "Región de Arica y Parinacota"
$ws.Range("P461").Value = 733
$ws.Range("Q461").Value = 15

# Row 462
$ws.Range("D462").Value = 44285
$ws.Range("H462").Value = "Cuatro cascos rojo"
$ws.Range("I462").Value = "Primera"
$ws.Range("J462").Value = 200
$ws.Range("K462").Value = 10000
$ws.Range("L462").Value = 10000
$ws.Range("M462").Value = 10000
$ws.Range("N462").Value = "`$/caja 15 kilos"
$ws.Range("O462").Value = "Región del Maule"
$ws.Range("P462").Value = 667
$ws.Range("Q462").Value = 15

# Row 463
$ws.Range("D463").Value = 44285
$ws.Range("H463").Value = "Cuatro cascos verde"
$ws.Range("I463").Value = "Primera"
$ws.Range("J463").Value = 200
$ws.Range("K463").Value = 7000
$ws.Range("L463").Value = 7000
$ws.Range("M463").Value = 7000
$ws.Range("N463").Value = "`$/caja 15 kilos"
$ws.Range("O463").Value = "Región del Maule"
$ws.Range("P463").Value = 467
$ws.Range("Q463").Value = 15

# Row 464
$ws.Range("D464").Value = 44498
$ws.Range("H464").Value = "Zafiro verde"
$ws.Range("I464").Value = "Primera"
$ws.Range("J464").Value = 300
$ws.Range("K464").Value = 35000
$ws.Range("L464").Value = 35000
$ws.Range("M464").Value = 35000
$ws.Range("N464").Value = "`$/caja 15 kilos"
$ws.Range("O464").Value = "Región de Arica y Parinacota"
$ws.Range("P464").Value = 2333
$ws.Range("Q464").Value = 15

# Row 465
$ws.Range("D465").Value = 44179
$ws.Range("H465").Value = "Cuatro cascos verde"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 200
$ws.Range("K465").Value = 12000
$ws.Range("L465").Value = 12000
$ws.Range("M465").Value = 12000
$ws.Range("N465").Value = "`$/caja 15 kilos"
$ws.Range("O465").Value = "Región del Maule"
$ws.Range("P465").Value = 800
$ws.Range("Q465").Value = 15

# Row 466
$ws.Range("D466").Value = 44418
$ws.Range("H466").Value = "Morrón rojo"
$ws.Range("I466").Value = "Primera"
$ws.Range("J466").Value = 200
$ws.Range("K466").Value = 38000
$ws.Range("L466").Value = 38000
$ws.Range("M466").Value = 38000
$ws.Range("N466").Value = "`$/caja 20 kilos"
$ws.Range("O466").Value = "Provincia del Elquí"
$ws.Range("P466").Value = 1900
$ws.Range("Q466").Value = 20

# Row 467
$ws.Range("D467").Value = 44418
$ws.Range("H467").Value = "Zafiro rojo"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 200
$ws.Range("K467").Value = 32000
$ws.Range("L467").Value = 32000
$ws.Range("M467").Value = 32000
$ws.Range("N467").Value = "`$/caja 15 kilos"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 2133
$ws.Range("Q467").Value = 15

# Row 468
$ws.Range("D468").Value = 44418
$ws.Range("H468").Value = "Zafiro verde"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 300
$ws.Range("K468").Value = 30000
$ws.Range("L468").Value = 30000
$ws.Range("M468").Value = 30000
$ws.Range("N468").Value = "`$/caja 15 kilos"
$ws.Range("O468").Value = "Región de Arica y Parinacota"
$ws.Range("P468").Value = 2000
$ws.Range("Q468").Value = 15

# Row 469
$ws.Range("D469").Value = 44595
$ws.Range("H469").Value = "Cuatro cascos verde"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 200
$ws.Range("K469").Value = 6000
$ws.Range("L469").Value = 6000
$ws.Range("M469").Value = 6000
$ws.Range("N469").Value = "`$/caja 15 kilos"
$ws.Range("O469").Value = "Región del Maule"
$ws.Range("P469").Value = 400
$ws.Range("Q469").Value = 15

# Row 470
$ws.Range("D470").Value = 44628
$ws.Range("H470").Value = "Cuatro cascos rojo"
$ws.Range("I470").Value = "Primera"
$ws.Range("J470").Value = 300
$ws.Range("K470").Value = 13000
$ws.Range("L470").Value = 13000
$ws.Range("M470").Value = 13000
$ws.Range("N470").Value = "`$/caja 15 kilos"
$ws.Range("O470").Value = "Región del Maule"
$ws.Range("P470").Value = 867
$ws.Range("Q470").Value = 15

# Row 471
$ws.Range("D471").Value = 44628
$ws.Range("H471").Value = "Cuatro cascos verde"
$ws.Range("I471").Value = "Primera"
$ws.Range("J471").Value = 300
$ws.Range("K471").Value = 8000
$ws.Range("L471").Value = 8000
$ws.Range("M471").Value = 8000
$ws.Range("N471").Value = "`$/caja 15 kilos"
$ws.Range("O471").Value = "Región del Maule"
$ws.Range("P471").Value = 533
$ws.Range("Q471").Value = 15

# Row 472
$ws.Range("D472").Value = 44335
$ws.Range("H472").Value = "Zafiro rojo"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 200
$ws.Range("K472").Value = 27000
$ws.Range("L472").Value = 27000
$ws.Range("M472").Value = 27000
$ws.Range("N472").Value = "`$/caja 15 kilos"
$ws.Range("O472").Value = "Región de Arica y Parinacota"
$ws.Range("P472").Value = 1800
$ws.Range("Q472").Value = 15

# Row 473
$ws.Range("D473").Value = 44335
$ws.Range("H473").Value = "Zafiro verde"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 400
$ws.Range("K473").Value = 13000
$ws.Range("L473").Value = 13000
$ws.Range("M473").Value = 13000
$ws.Range("N473").Value = "`$/caja 15 kilos"
$ws.Range("O473").Value = "Región de Arica y Parinacota"
$ws.Range("P473").Value = 867
$ws.Range("Q473").Value = 15

# Row 474
$ws.Range("A474").Value = 5
$ws.Range("B474").Value = "Macroferia Regional de Talca"
$ws.Range("C474").Value = "Maule"
$ws.Range("D474").Value = 44552
$ws.Range("E474").Value = 7
$ws.Range("F474").Value = 100112002
$ws.Range("G474").Value = "Pimiento"
$ws.Range("H474").Value = "Cuatro cascos verde"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 200
$ws.Range("K474").Value = 12000
$ws.Range("L474").Value = 12000
$ws.Range("M474").Value = 12000
$ws.Range("N474").Value = "`$/caja 15 kilos"
$ws.Range("O474").Value = "Región del Maule"
$ws.Range("P474").Value = 800
$ws.Range("Q474").Value = 15
$ws.Range("R474").Value = "Hortaliza"

# Row 475
$ws.Range("A475").Value = 5
$ws.Range("B475").Value = "Macroferia Regional de Talca"
$ws.Range("C475").Value = "Maule"
$ws.Range("D475").Value = 44544
$ws.Range("E475").Value = 7
$ws.Range("F475").Value = 100112002
$ws.Range("G475").Value = "Pimiento"
$ws.Range("H475").Value = "Cuatro cascos verde"
$ws.Range("I475").Value = "Primera"
$ws.Range("J475").Value = 300
$ws.Range("K475").Value = 12000
$ws.Range("L475").Value = 12000
$ws.Range("M475").Value = 12000
$ws.Range("N475").Value = "`$/caja 15 kilos"
$ws.Range("O475").Value = "Región del Maule"
$ws.Range("P475").Value = 800
$ws.Range("Q475").Value = 15
$ws.Range("R475").Value = "Hortaliza"

$dfmt = $ws.Range("D473").NumberFormat()
$ws.Range("D474").NumberFormat = $dfmt
$ws.Range("D475").NumberFormat = $dfmt

Write-Host "Done"